$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 4, shifting the existing rows 4-40 down to 5-41.
# Excel's Insert() carries the formatting of the row above down with the
# shifted rows, so the date-formatted column D keeps its number format.
$ws.Rows("4:4").Insert()

# Populate the newly inserted row 4 with the new weekly record.
$ws.Range("A4").Value2 = 5
$ws.Range("B4").Value2 = "Macroferia Regional de Talca"
$ws.Range("C4").Value2 = "Maule"
$ws.Range("D4").Value2 = 44685
$ws.Range("E4").Value2 = 7
$ws.Range("F4").Value2 = "Fruta"
$ws.Range("G4").Value2 = 100107
$ws.Range("H4").Value2 = "Otros"
$ws.Range("I4").Value2 = 100107001
$ws.Range("J4").Value2 = "Caqui"
$ws.Range("K4").Value2 = "Mankaki"
$ws.Range("L4").Value2 = "Primera"
$ws.Range("M4").Value2 = 180
$ws.Range("N4").Value2 = 23000
$ws.Range("O4").Value2 = 23000
$ws.Range("P4").Value2 = 23000
$ws.Range("Q4").Value2 = "`$/caja 18 kilos granel"
$ws.Range("R4").Value2 = "Región del Maule"
$ws.Range("S4").Value2 = 1278
$ws.Range("T4").Value2 = 18
